# Se arregla recibos de pago
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E became much narrower (stored OOXML width 28 -> 7).
# ColumnWidth = 6.083333333333333 reliably round-trips to a stored width of 7.
$ws.Columns.Item(5).ColumnWidth = 6.083333333333333

# ---- Block 1 (rows 4-22) ----
$ws.Range("A4").Value = "JAIME MARTINEZ"
$ws.Range("D4").Value = "0-441974958290"
$ws.Range("G4").Value = 130.79

$ws.Range("C8").Value = 4.51
$ws.Range("G8").Value = 130.79

$ws.Range("C9").Value = 36.08
$ws.Range("G9").Value = 130.79

$ws.Range("C10").Value = 29

$ws.Range("G12").Value = 130.79

$ws.Range("B13").Value = 29
$ws.Range("C13").Value = 130.79

$ws.Range("D15").Value = "BANCO GENERAL"

$ws.Range("B17").Value = 29
$ws.Range("C17").Value = 130.79

$ws.Range("B20").Value = 130.79
$ws.Range("C20").Value = "BANCO GENERAL"
$ws.Range("G20").Value = 130.79

$ws.Range("A22").Value = "JAIME MARTINEZ"

# ---- Block 2 (rows 30-48) ----
$ws.Range("A30").Value = "JAIME MARTINEZ"
$ws.Range("D30").Value = "0-441974958290"
$ws.Range("G30").Value = 130.79

$ws.Range("C34").Value = 4.51
$ws.Range("G34").Value = 130.79

$ws.Range("C35").Value = 36.08
$ws.Range("G35").Value = 130.79

$ws.Range("C36").Value = 29

$ws.Range("G38").Value = 130.79

$ws.Range("B39").Value = 29
$ws.Range("C39").Value = 130.79

$ws.Range("D41").Value = "BANCO GENERAL"

$ws.Range("B43").Value = 29
$ws.Range("C43").Value = 130.79

$ws.Range("B46").Value = 130.79
$ws.Range("C46").Value = "BANCO GENERAL"
$ws.Range("G46").Value = 130.79

$ws.Range("A48").Value = "JAIME MARTINEZ"
